# LOQ4268.xlsx update
# - Insert a new row at 13 (shifts the "Docentes responsáveis" value down and
#   everything below it down by one row), then fix up the B/C (value /
#   highlighted-value) columns so the sheet matches the refreshed course
#   description: new "Objetivos" text, new "Programa resumido" summary, the
#   full "Programa" syllabus text, the "Método"/"Critério"/"Norma de
#   recuperação" texts in their corrected rows, and a new "Bibliografia" entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 13-23 down to 14-24 (creates a fresh, unformatted row 13).
$ws.Rows.Item(13).Insert()

# The freshly-inserted row 13 has no real B/C formatting yet (and an
# empty/unused A13 cell) — borrow the normal-value/red-value cell formats
# from a row that still has them, then drop the unused A13 cell.
$ws.Range("B10").Copy()
$ws.Range("B13").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("C10").Copy()
$ws.Range("C13").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false

# Row 10 "Objetivos:" — replace the (mis-copied) value with the real objectives text.
$ws.Range("B10").Value = "Oferecer ao aluno uma visão geral dos conceitos fundamentais e das fases do processo de desenvolvimento e Implementação de Sistemas de Informação no sentido de capacitá-lo analisa e projetar tais sistemas"
$ws.Range("C10").Value = "Oferecer ao aluno uma visão geral dos conceitos fundamentais e das fases do processo de desenvolvimento e Implementação de Sistemas de Informação no sentido de capacitá-lo analisa e projetar tais sistemas"

# Row 13 (new row, under "Docentes responsáveis:") — professor name.
$ws.Range("B13").Value = "11079086 - Herlandí de Souza Andrade"
$ws.Range("C13").Value = "11079086 - Herlandí de Souza Andrade"

# Row 14 "Programa resumido:" — short syllabus summary.
$ws.Range("B14").Value = "Sistemas de Informação. Projeto de Sistemas de Informação. Tecnologia de Informação. Processo de Desenvolvimento de Sistema de Informação."
$ws.Range("C14").Value = "Sistemas de Informação. Projeto de Sistemas de Informação. Tecnologia de Informação. Processo de Desenvolvimento de Sistema de Informação."

# Row 16 "Programa:" — full syllabus text.
$ws.Range("B16").Value = "1. Sistemas de Informação1.1. Sistemas de Processamento de Informações;1.2. Sistemas de Informações Gerenciais;1.3. Sistema de Apoio à Decisão;1.4. Sistemas de Informação no Comércio Eletrônico;1.5. Sistemas de Informação em Cadeia de Suprimentos;1.6. Sistemas inteligentes nos negócios;1.7. Sistemas estratégicos. 2. Projeto de Sistemas de Informação.2.1. Especificação das Saídas;2.2. Especificação dos Arquivos;2.3. Especificação das Entradas;2.4. Especificação do Processamento.3. Tecnologia de Informação.3.1. Evolução da Computação;3.2. Recursos Computacionais.4. Processo de Desenvolvimento de Sistemas de Informação.4.1. Definição do Negócio;4.2. Identificação do Problema e/ou Oportunidades;4.3. Seleção do Sistema de Informação;4.4. Implementação do Sistema de Informação;4.5. Avaliação da Eficácia do Sistema de Informação;"
$ws.Range("C16").Value = "1. Sistemas de Informação1.1. Sistemas de Processamento de Informações;1.2. Sistemas de Informações Gerenciais;1.3. Sistema de Apoio à Decisão;1.4. Sistemas de Informação no Comércio Eletrônico;1.5. Sistemas de Informação em Cadeia de Suprimentos;1.6. Sistemas inteligentes nos negócios;1.7. Sistemas estratégicos. 2. Projeto de Sistemas de Informação.2.1. Especificação das Saídas;2.2. Especificação dos Arquivos;2.3. Especificação das Entradas;2.4. Especificação do Processamento.3. Tecnologia de Informação.3.1. Evolução da Computação;3.2. Recursos Computacionais.4. Processo de Desenvolvimento de Sistemas de Informação.4.1. Definição do Negócio;4.2. Identificação do Problema e/ou Oportunidades;4.3. Seleção do Sistema de Informação;4.4. Implementação do Sistema de Informação;4.5. Avaliação da Eficácia do Sistema de Informação;"

# Row 19 "Método:" — teaching method text (now correctly under "Método").
$ws.Range("B19").Value = "Aulas expositivas teóricas, aulas práticas, aulas de exercícios."
$ws.Range("C19").Value = "Aulas expositivas teóricas, aulas práticas, aulas de exercícios."

# Row 20 "Critério:" — evaluation criteria text.
$ws.Range("B20").Value = "Média Aritmética das atividades avaliativas realizadas."
$ws.Range("C20").Value = "Média Aritmética das atividades avaliativas realizadas."

# Row 21 "Norma de recuperação:" — recovery grading rule text.
$ws.Range("B21").Value = "Média aritmética da nota final obtida pelo aluno durante o semestre e da nota obtida na Prova de Recuperação."
$ws.Range("C21").Value = "Média aritmética da nota final obtida pelo aluno durante o semestre e da nota obtida na Prova de Recuperação."

# Row 22 "Bibliografia:" — bibliography text.
$ws.Range("B22").Value = "HAL R. VARIAN, H. R.; FARRELL, J., SHAPIRO, C. The economics of information technology: an introduction. Cambridge University Press, 2004LAURINDO, F. J. B.; ROTONDARO, R. G. orgs. Gestão integrada de processos e da tecnologia da informação. São Paulo:Atlas, 2006.LAURINDO, F.J.B. Tecnologia da Informação: Eficácia nas Organizações. São Paulo, Editora Futura, 2002.STAIR, R.M., Princípios de Sistema de Informação: Uma Abordagem Gerencial, Rio de Janeiro, LTC, 1998.TURBAN, E. et al. Information Technology for Management: Transforming Organizations in the Digital Economy. 7th edition, Wiley, 2009.TURBAN, E., RAIANER JR, K., POTTER, R. E., Administração de Tecnologia da Informação: Teoria e Prática”, São Paulo, Editora Campus, 2003."
$ws.Range("C22").Value = "HAL R. VARIAN, H. R.; FARRELL, J., SHAPIRO, C. The economics of information technology: an introduction. Cambridge University Press, 2004LAURINDO, F. J. B.; ROTONDARO, R. G. orgs. Gestão integrada de processos e da tecnologia da informação. São Paulo:Atlas, 2006.LAURINDO, F.J.B. Tecnologia da Informação: Eficácia nas Organizações. São Paulo, Editora Futura, 2002.STAIR, R.M., Princípios de Sistema de Informação: Uma Abordagem Gerencial, Rio de Janeiro, LTC, 1998.TURBAN, E. et al. Information Technology for Management: Transforming Organizations in the Digital Economy. 7th edition, Wiley, 2009.TURBAN, E., RAIANER JR, K., POTTER, R. E., Administração de Tecnologia da Informação: Teoria e Prática”, São Paulo, Editora Campus, 2003."

# Row heights: rows 17 and 22 grow a tall "customHeight" (Syllabus / Bibliografia
# now hold long wrapped text), row 21 shrinks back to the normal 60pt height.
$ws.Rows.Item(17).RowHeight = 120
$ws.Rows.Item(21).RowHeight = 60
$ws.Rows.Item(22).RowHeight = 120

# Row 13 has no "label" column — only the professor name in B/C — so drop
# the leftover, unused A13 cell entirely.
$ws.Range("A13").Clear()
